$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record (Fecha 2021-08-30 / Volumen 3000) was inserted at the
# top of the data block (row 139), pushing every existing record from row
# 139 down one row (139->140, 140->141, ... 233->234).
$ws.Rows.Item(139).Insert()

# After the insert, row 140 holds what used to live in row 139 (all columns
# shifted down intact). Seed the freshly-inserted, still-empty row 139 with
# that same record...
$src = $ws.Range("A140:R140")
$dst = $ws.Range("A139:R139")
$dst.Value = $src.Value()

# ...then overwrite just the Fecha (D) and Volumen (J) columns with the new
# values for this newly-added record.
$ws.Range("D139").Value = 44438
$ws.Range("J139").Value = 3000
